$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 with the new combined/stringified card data
$ws.Range("A2").Value = "('Soul of Ravnica', ['{4}{U}{U}', 'Creature " + [char]8212 + " Avatar', 'Flying', '{5}{U}{U}: Draw a card for each color among permanents you control.', '{5}{U}{U}, Exile Soul of Ravnica from your graveyard: Draw a card for each color among permanents you control.', '6/6'])"
$ws.Range("A3").Value = "('Soul of Zendikar', ['{4}{G}{G}', 'Creature " + [char]8212 + " Avatar', 'Reach', '{3}{G}{G}: Create a 3/3 green Beast creature token.', '{3}{G}{G}, Exile Soul of Zendikar from your graveyard: Create a 3/3 green Beast creature token.', '6/6'])"

# Clear the now-unused rows 4 through 15
$ws.Range("A4:A15").Clear()
